# Generate Report for Handback
# Updates the "Correspond Handback DateTime" column (G) for each locale
# sheet with a refreshed report-generation timestamp, leaving any rows
# that already carry a distinct handback datetime untouched.

$wb = $excel.ActiveWorkbook

# zh-cn: every data row (2-102) shares one handback datetime value.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("G2:G102").Value = "2016-02-22 08:50:29"

# de-de: most data rows share one handback datetime value, but rows
# 7, 9, 11, 80 and 101 already record their own distinct datetimes
# and must be left alone.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("G2:G6").Value = "2016-02-22 08:50:40"
$wsDeDe.Range("G8:G8").Value = "2016-02-22 08:50:40"
$wsDeDe.Range("G10:G10").Value = "2016-02-22 08:50:40"
$wsDeDe.Range("G12:G79").Value = "2016-02-22 08:50:40"
$wsDeDe.Range("G81:G100").Value = "2016-02-22 08:50:40"
$wsDeDe.Range("G102:G105").Value = "2016-02-22 08:50:40"

# ja-jp: most data rows share one handback datetime value, but row 78
# already records its own distinct datetime and must be left alone.
$wsJaJp = $wb.Worksheets.Item("ja-jp")
$wsJaJp.Range("G2:G77").Value = "2016-02-22 08:50:50"
$wsJaJp.Range("G79:G106").Value = "2016-02-22 08:50:50"

# zh-tw: every data row (2-105) shares one handback datetime value.
$wsZhTw = $wb.Worksheets.Item("zh-tw")
$wsZhTw.Range("G2:G105").Value = "2016-02-22 08:51:01"
